$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: confirmation text for the already-collected amount (matches B1)
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "1000"

# Row 2: confirmation text that does NOT match B2 (mismatched entry kept as-is)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "999"

# Row 3: new collection entry
$ws.Range("B3").Value = 123
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "123"

# Row 4: new collection entry
$ws.Range("B4").Value = 135
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "135"

# Row 5: new collection entry
$ws.Range("B5").Value = 133
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "133"

# Leave the selection where Excel would land after entering the last row
$ws.Range("B6").Select() | Out-Null
